# Regenerate merged AHB files
# - Rename the "_old" / "_new" header suffixes to the new release tags
#   (FV2410 / FV2504).
# - Wrap the data range in an Excel Table ("Table1") with an AutoFilter.
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row ---------------------------------------------
# Columns A:J carried the "_old" suffix -> becomes "_FV2410"
$fv2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $fv2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410[$i]
}

# Column K ("diff") is unchanged.

# Columns L:U carried the "_new" suffix -> becomes "_FV2504"
$fv2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)
for ($i = 0; $i -lt $fv2504.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504[$i]
}

# --- 2. Freeze the header row -------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel Table with an AutoFilter ------
$dataRange = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

Write-Host "AHB merge regenerated: headers renamed, Table1 created, header row frozen."
